$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume/hour data refresh (GitHub Actions run).
# Values are written as text to match the sheet's existing inline-string cell type,
# so each cell is briefly formatted as Text ("@") before the write, then the
# temporary number-format override is cleared again.
$updates = @(
    @('D2', '302.34'),
    @('E2', '-0.39%'),
    @('G2', '9'),
    @('D3', '37.42'),
    @('E3', '7.05%'),
    @('G3', '9'),
    @('D4', '5.005'),
    @('E4', '-4.19%'),
    @('G4', '9'),
    @('E5', '0.26%'),
    @('G5', '9'),
    @('D6', '2.193'),
    @('E6', '-8.15%'),
    @('G6', '9'),
    @('D7', '8.042'),
    @('E7', '-0.11%'),
    @('G7', '9'),
    @('D8', '4.040'),
    @('E8', '1.58%'),
    @('G8', '9'),
    @('D9', '0.9117'),
    @('E9', '-2.52%'),
    @('G9', '9'),
    @('D10', '0.09688'),
    @('E10', '-3.59%'),
    @('G10', '9'),
    @('D11', '0.1883'),
    @('E11', '2.35%'),
    @('G11', '9'),
    @('D12', '0.08579'),
    @('E12', '-1.78%'),
    @('G12', '9'),
    @('D13', '0.03520'),
    @('E13', '6.05%'),
    @('G13', '9'),
    @('D14', '0.09979'),
    @('E14', '0.68%'),
    @('G14', '9'),
    @('D15', '0.001477'),
    @('E15', '-0.63%'),
    @('G15', '9'),
    @('D16', '0.005658'),
    @('E16', '-0.03%'),
    @('G16', '9'),
    @('E17', '0.03%'),
    @('G17', '9'),
    @('D18', '2.077'),
    @('E18', '-3.43%'),
    @('G18', '9'),
    @('E19', '2.55%'),
    @('G19', '9'),
    @('D20', '0.1302'),
    @('E20', '-0.03%'),
    @('G20', '9'),
    @('D21', '4.762'),
    @('E21', '10.27%'),
    @('G21', '9'),
    @('D22', '0.2206'),
    @('E22', '-1.01%'),
    @('G22', '9'),
    @('D23', '0.04628'),
    @('E23', '1.19%'),
    @('G23', '9'),
    @('E24', '1.10%'),
    @('G24', '9'),
    @('D25', '0.004804'),
    @('E25', '8.39%'),
    @('G25', '9'),
    @('E26', '-7.66%'),
    @('G26', '9'),
    @('E27', '28.45%'),
    @('G27', '9'),
    @('G28', '9'),
    @('G29', '9'),
    @('G30', '9'),
    @('G31', '9'),
    @('G32', '9'),
    @('G33', '9'),
    @('G34', '9'),
    @('G35', '9'),
    @('G36', '9'),
    @('G37', '9'),
    @('G38', '9'),
    @('D39', '0.01761'),
    @('E39', '-0.26%'),
    @('G39', '9'),
    @('D40', '0.04726'),
    @('E40', '-1.75%'),
    @('G40', '9'),
    @('D41', '0.008058'),
    @('E41', '3.77%'),
    @('G41', '9'),
    @('E42', '-1.15%'),
    @('G42', '9'),
    @('D43', '0.007537'),
    @('E43', '5.69%'),
    @('G43', '9'),
    @('D44', '0.002181'),
    @('E44', '-1.02%'),
    @('G44', '9'),
    @('D45', '0.01042'),
    @('G45', '9'),
    @('D46', '0.00006060'),
    @('E46', '4.06%'),
    @('G46', '9'),
    @('D47', '0.00000000750'),
    @('E47', '0.02%'),
    @('G47', '9'),
    @('D48', '8.671'),
    @('E48', '217.34%'),
    @('G48', '9'),
    @('G49', '9'),
    @('E50', '0.02%'),
    @('G50', '9'),
    @('E51', '0.02%'),
    @('G51', '9'),
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newValue = $u[1]
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.ClearFormats()
}
